# Insert a new data row at row 765 (pushing the existing rows 765-817
# down to 766-818) and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("765:765").Insert()

$ws.Range("A765").Value = 4
$ws.Range("B765").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C765").Value = "Los Lagos"
$ws.Range("D765").Value = 45265
$ws.Range("E765").Value = 10
$ws.Range("F765").Value = 100114001
$ws.Range("G765").Value = "Papa"
$ws.Range("H765").Value = "Patagonia"
$ws.Range("I765").Value = "1a nueva(o)"
$ws.Range("J765").Value = 600
$ws.Range("K765").Value = 25000
$ws.Range("L765").Value = 26000
$ws.Range("M765").Value = 25500
$ws.Range("N765").Value = "$/saco 25 kilos"
$ws.Range("O765").Value = "Región de La Araucanía"
$ws.Range("P765").Value = 1020
$ws.Range("Q765").Value = 25
$ws.Range("R765").Value = "Hortaliza"
